# Update countries & provincias Spain
# Refresh of the COVID-19 "Pais" dashboard data: updated totals for several
# countries (which, since the sheet is kept sorted descending by "Casos
# totales", causes Barein to overtake Azerbaiyan and Jordania to overtake
# Reunion/Taiwan), plus a refreshed "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 19:52"

# Estados Unidos (row 4): Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$ws.Range("B4").Value = 549633
$ws.Range("C4").Value = 16754
$ws.Range("D4").Value = 31120
$ws.Range("E4").Value = 496872
$ws.Range("G4").Value = 1064
$ws.Range("H4").Value = 21641

# Alemania (row 8)
$ws.Range("D8").Value = 60300
$ws.Range("E8").Value = 63448

# Suiza (row 15)
$ws.Range("D15").Value = 12700
$ws.Range("E15").Value = 11601

# Pakistan (row 36)
$ws.Range("B36").Value = 5230
$ws.Range("C36").Value = 219
$ws.Range("E36").Value = 4111
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 91

# Marruecos (row 61)
$ws.Range("B61").Value = 1661
$ws.Range("C61").Value = 116
$ws.Range("D61").Value = 177
$ws.Range("E61").Value = 1366
$ws.Range("G61").Value = 7
$ws.Range("H61").Value = 118

# Barein overtakes Azerbaiyan: row 69 becomes Barein (updated data),
# row 70 becomes Azerbaiyan (its previous, unchanged data)
$ws.Range("A69").Value = "Barein"
$ws.Range("B69").Value = 1136
$ws.Range("C69").Value = 96
$ws.Range("D69").Value = 558
$ws.Range("E69").Value = 572
$ws.Range("F69").Value = 4
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 6

$ws.Range("A70").Value = "Azerbaiyan"
$ws.Range("B70").Value = 1098
$ws.Range("C70").Value = 40
$ws.Range("D70").Value = 250
$ws.Range("E70").Value = 837
$ws.Range("F70").Value = 27
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 11

# Kazajistan (row 75)
$ws.Range("B75").Value = 951
$ws.Range("C75").Value = 86
$ws.Range("E75").Value = 842

# Libano (row 87)
$ws.Range("D87").Value = 80
$ws.Range("E87").Value = 530

# Jordania overtakes Reunion and Taiwan: row 99 becomes Jordania (updated
# data), row 100 becomes Reunion (its previous, unchanged data), row 101
# becomes Taiwan (its previous, unchanged data)
$ws.Range("A99").Value = "Jordania"
$ws.Range("B99").Value = 389
$ws.Range("C99").Value = 8
$ws.Range("D99").Value = 201
$ws.Range("E99").Value = 181
$ws.Range("F99").Value = 5
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 7

$ws.Range("A100").Value = "Reunion"
$ws.Range("B100").Value = 388
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 40
$ws.Range("E100").Value = 348
$ws.Range("F100").Value = 3
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0

$ws.Range("A101").Value = "Taiwan"
$ws.Range("B101").Value = 388
$ws.Range("C101").Value = 3
$ws.Range("D101").Value = 109
$ws.Range("E101").Value = 273
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 6

# Georgia (row 112)
$ws.Range("B112").Value = 257
$ws.Range("C112").Value = 15
$ws.Range("E112").Value = 187
